# Updates cryptos list values (price/volume columns, plus a few row-order swaps)
# as captured by the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.086.05"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.427.61"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.34%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.03"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.59"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.28%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.421.84"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.39%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.55%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.94"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.69%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.561"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.47%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "43.98"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.40%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000268"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.982.66"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.09"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("E16").ClearFormats()
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("B17").ClearFormats()
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C17").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.449.34"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("E17").ClearFormats()
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("B18").ClearFormats()
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("C18").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.211.44"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("E18").ClearFormats()
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("B19").ClearFormats()
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C19").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "579.03"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.48%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.120"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.03"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.840"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.80"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.21%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "95.30"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.02"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.43%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.63"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.44%  "
$ws.Range("E26").ClearFormats()
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.40"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.15%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.52"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.90%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.51"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.65%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.76"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.12%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.23"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.09%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.75"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -7.37%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.52"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.59%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "574.53"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -9.63%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.40"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.68%  "
$ws.Range("E36").ClearFormats()
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("B37").ClearFormats()
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C37").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("E37").ClearFormats()
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "VeChain"
$ws.Range("B38").ClearFormats()
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C38").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0467"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.23%  "
$ws.Range("E38").ClearFormats()
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Hedera"
$ws.Range("B39").ClearFormats()
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C39").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0949"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.42%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "55.92"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.140"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.10"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -13.78%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.223.28"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.26%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0676"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -8.63%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "30.91"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.92%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.292"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.85%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.75"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.12%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.37"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -6.98%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.125"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.41%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.86"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("E50").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.02%  "
$ws.Range("E51").ClearFormats()
